$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for the time-variable covariate columns
$ws.Range("G1").Value = "WT"
$ws.Range("H1").Value = "CRCL"
$ws.Range("I1").Value = "DIAL"

# Fill the new covariate columns with the same "missing data" marker already
# used in the CONC column (D) for rows 2-6
$ws.Range("G2:I6").Value = "."

# Touch row 7 so Excel materializes an (empty) row element for it, matching
# the row height already used by the other blank formatted rows below it
$ws.Range("A7").RowHeight = 13.8
$ws.Range("A7").Select()

# Excel keeps a handful of formatted-but-empty rows trailing the used range;
# extend that trailing block by four more rows
$ws.Range("A1048571:A1048574").RowHeight = 12.8
